$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.734267234802246
$ws.Range("B1").Value = 2.463452577590942
$ws.Range("C1").Value = 4.84526252746582
$ws.Range("D1").Value = 4.447070598602295
$ws.Range("E1").Value = 1.261392831802368
